# Applies the cryptos list price/volume refresh described in the commit
# 'Updated cryptos list on Mon Oct 28 05:52:01 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) are plain text cells in this sheet.
# Some new price strings (e.g. '587.60') look numeric, so force the cell's
# number format to Text before assigning them, to avoid Excel auto-converting
# them into numbers (which would drop formatting such as trailing zeros).

$ws.Range('D2').Value = '67.807.58'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.484.57'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.60'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.19'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('E9').Value = '  +3.74%  '
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '2.937.87'
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.28'
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').Value = '67.681.69'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '2.485.62'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.80'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.37'
$ws.Range('E19').Value = '  -2.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '346.23'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.10'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.67'
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('E25').Value = '  -7.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.83'
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('D27').Value = '2.611.45'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').Value = '0.0₃0887'
$ws.Range('E29').Value = '  -2.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '498.05'
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.71'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '164.53'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.64'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('E40').Value = '  -2.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.72'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.324'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.76'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.32'
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.510'
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('E48').Value = '  -4.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0734'
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('E50').Value = '  -1.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.575'
$ws.Range('E51').Value = '  -1.50%  '
